# Update the "Puntuacion Media" (average rating) column so that it now
# reflects user comments instead of the previous static value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 5
$ws.Range("E3").Value = 5

# Leave the selection where the author left it when they saved the file.
$ws.Range("H7").Select()
